$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 10; $r++) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $bVal = $bCell.Value2
    $cVal = $cCell.Value2
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2

    $bCell.Value2 = $cVal
    $cCell.Value2 = $bVal
    $dCell.Value2 = $eVal
    $eCell.Value2 = $dVal
}
